$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.177.83"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "3.105.72"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "574.25"
$c.ClearFormats()
$ws.Range("E5").Value = "  -0.82%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "178.42"
$c.ClearFormats()
$ws.Range("E6").Value = "  +4.12%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.106.01"
$ws.Range("E8").Value = "  +0.63%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.512"
$c.ClearFormats()
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("E13").Value = "  -1.79%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "36.23"
$c.ClearFormats()
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "3.625.75"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "67.118.79"
$ws.Range("E17").Value = "  +0.14%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "7.03"
$c.ClearFormats()
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").Value = "3.108.72"
$ws.Range("E19").Value = "  +0.61%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "16.70"
$c.ClearFormats()
$ws.Range("E20").Value = "  +1.15%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "491.35"
$c.ClearFormats()
$ws.Range("E21").Value = "  +0.91%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "7.72"
$c.ClearFormats()
$ws.Range("E22").Value = "  -0.39%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.687"
$c.ClearFormats()
$ws.Range("E23").Value = "  -1.12%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "83.49"
$c.ClearFormats()
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("E25").Value = "  +0.74%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "12.59"
$c.ClearFormats()
$ws.Range("E26").Value = "  -2.61%  "
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("E28").Value = "  +0.06%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.88"
$c.ClearFormats()
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("E30").Value = "  -0.13%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.60"
$c.ClearFormats()
$ws.Range("E31").Value = "  -2.00%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "28.15"
$c.ClearFormats()
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("D34").Value = "0.0₃0941"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.949"
$c.ClearFormats()
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("B37").Value = "Arweave"
$ws.Range("C37").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "47.29"
$c.ClearFormats()
$ws.Range("E37").Value = "  +1.92%  "
$ws.Range("E38").Value = "  -2.94%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.312"
$c.ClearFormats()
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("E40").Value = "  +0.82%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "49.19"
$c.ClearFormats()
$ws.Range("E41").Value = "  -1.07%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "8.25"
$c.ClearFormats()
$ws.Range("E43").Value = "  -1.82%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.71"
$c.ClearFormats()
$ws.Range("E44").Value = "  +5.57%  "
$ws.Range("D45").Value = "2.802.69"
$ws.Range("E45").Value = "  +0.59%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "371.21"
$c.ClearFormats()
$ws.Range("E46").Value = "  -2.15%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0345"
$c.ClearFormats()
$ws.Range("E47").Value = "  -0.86%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "135.57"
$c.ClearFormats()
$ws.Range("E48").Value = "  +0.48%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "25.55"
$c.ClearFormats()
$ws.Range("E50").Value = "  +3.27%  "
$ws.Range("E51").Value = "  +5.11%  "
